$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Append the new paragraphs after the final "Data Boot Camp..." line.
# All text is inserted first (plain), then formatting (bold run,
# first-line indents, the relocated _GoBack bookmark) is applied
# afterwards via Find so that paragraph marks created along the way
# don't pick up stray character formatting.
# ------------------------------------------------------------------

# Paragraph: empty (NoSpacing)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "No Spacing"

# Paragraph: "#Weather Analysis:" (NoSpacing)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "No Spacing"
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("#Weather Analysis:")

# Paragraph: "The plot of latitude vs weather indicates ... latitude. " (NoSpacing, firstLine indent)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "No Spacing"
$d.Paragraphs.Last.Format.FirstLineIndent = 36
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("The plot of latitude vs weather indicates that the temperature is highest near the Tropic of cancer with latitude around 23.5 deg N. T The maximum temperatures in the northern hemisphere falls considerably moving north of the 23.5 deg latitude and moving south from the equator at 0 deg. There is also a wide variation in the maximum temperatures around the world based on the latitude. ")

# Paragraph: empty (NoSpacing, firstLine indent)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "No Spacing"
$d.Paragraphs.Last.Format.FirstLineIndent = 36

# Paragraph: single space (NoSpacing)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "No Spacing"
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter(" ")

# ------------------------------------------------------------------
# Bold "latitude vs weather" inside the paragraph we just typed.
# ------------------------------------------------------------------
$fr = $d.Content
$fr.Find.Execute("latitude vs weather", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Bold = 1

# ------------------------------------------------------------------
# Re-create the _GoBack bookmark so it wraps the "#Weather Analysis:"
# heading through the end of the "...based on the latitude. "
# paragraph. Word keeps only a single _GoBack bookmark at a time, so
# this both places the new one and removes the stale one that used to
# sit right after "Interpret the results."
# ------------------------------------------------------------------
$bmStartR = $d.Content
$bmStartR.Find.Execute("#Weather Analysis:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmStart = $bmStartR.Start

$bmEndR = $d.Content
$bmEndR.Find.Execute("based on the latitude. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmEnd = $bmEndR.End

$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
